$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 3000
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 3000
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 9000
$ws.Range("M69").ClearContents()
$ws.Range("N69").Value = -10748
$ws.Range("H72").Value = 3000
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 3000
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 27000
$ws.Range("M72").ClearContents()
$ws.Range("N72").Value = -35736
$ws.Range("H96").Value = 1640.6111
$ws.Range("I96").Value = 1372.8889
$ws.Range("J96").Value = 1908.3334
$ws.Range("K96").Value = 4118.6667
$ws.Range("L96").Value = 5725.0002
$ws.Range("M96").Value = -2745.6667
$ws.Range("N96").Value = -8471.0002
$ws.Range("H100").Value = 2522.6924
$ws.Range("I100").Value = 2368.4614
$ws.Range("J100").Value = 2676.923
$ws.Range("K100").Value = 2368.4614
$ws.Range("L100").Value = 2676.923
$ws.Range("M100").Value = -1827.4614
$ws.Range("N100").Value = -3758.923
$ws.Range("H137").Value = 14667.333
$ws.Range("I137").Value = 14667.333
$ws.Range("K137").Value = 44001.999
$ws.Range("M137").Value = -41451.999
$ws.Range("H138").Value = 249518.6
$ws.Range("I138").Value = 4706.625
$ws.Range("J138").Value = 303921.25
$ws.Range("K138").Value = 14119.875
$ws.Range("L138").Value = 911763.75
$ws.Range("M138").Value = -8979.875
$ws.Range("N138").Value = -922043.75

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3000.2666
$ws.Range("I61").Value = 2079
$ws.Range("J61").Value = 4842.8
$ws.Range("K61").Value = 2079
$ws.Range("L61").Value = 4842.8
$ws.Range("M61").Value = -1867
$ws.Range("N61").Value = -5266.8
$ws.Range("H74").Value = 1660.6666
$ws.Range("I74").Value = 1151
$ws.Range("J74").Value = 2680
$ws.Range("K74").Value = 1151
$ws.Range("L74").Value = 2680
$ws.Range("M74").Value = -277
$ws.Range("N74").Value = -4428
$ws.Range("H77").Value = 1660.6666
$ws.Range("I77").Value = 1151
$ws.Range("J77").Value = 2680
$ws.Range("K77").Value = 5755
$ws.Range("L77").Value = 13400
$ws.Range("M77").Value = -1387
$ws.Range("N77").Value = -22136
$ws.Range("H110").Value = 1875.75
$ws.Range("I110").Value = 1790.9
$ws.Range("J110").Value = 2300
$ws.Range("K110").Value = 1790.9
$ws.Range("L110").Value = 2300
$ws.Range("M110").Value = 254.0999999999999
$ws.Range("N110").Value = -6390
$ws.Range("H132").Value = 4323.4707
$ws.Range("I132").Value = 3475.762
$ws.Range("K132").Value = 10427.286
$ws.Range("M132").Value = -7897.286
$ws.Range("H136").Value = 3000.2666
$ws.Range("I136").Value = 2079
$ws.Range("J136").Value = 4842.8
$ws.Range("K136").Value = 6237
$ws.Range("L136").Value = 14528.4
$ws.Range("M136").Value = -3687
$ws.Range("N136").Value = -19628.4

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 7354914.5
$ws.Range("I105").Value = 7354914.5
$ws.Range("K105").Value = 7354914.5
$ws.Range("M105").Value = -7353167.5
$ws.Range("H107").Value = 1628.5714
$ws.Range("I107").Value = 966.6667
$ws.Range("J107").Value = 2125
$ws.Range("K107").Value = 966.6667
$ws.Range("L107").Value = 2125
$ws.Range("M107").Value = 953.3333
$ws.Range("N107").Value = -5965
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()
$ws.Range("H134").Value = 4090.7222
$ws.Range("I134").Value = 4025.889
$ws.Range("J134").Value = 4155.5557
$ws.Range("K134").Value = 12077.667
$ws.Range("L134").Value = 12466.6671
$ws.Range("M134").Value = -9542.667000000001
$ws.Range("N134").Value = -17536.6671

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4593.1562
$ws.Range("I31").Value = 1072.9375
$ws.Range("K31").Value = 1072.9375
$ws.Range("M31").Value = -777.9375
$ws.Range("H34").Value = 4593.1562
$ws.Range("I34").Value = 1072.9375
$ws.Range("K34").Value = 1072.9375
$ws.Range("M34").Value = -870.9375
$ws.Range("H58").Value = 1364.7894
$ws.Range("I58").Value = 1115.875
$ws.Range("J58").Value = 1545.8182
$ws.Range("K58").Value = 1115.875
$ws.Range("L58").Value = 1545.8182
$ws.Range("M58").Value = -912.875
$ws.Range("N58").Value = -1951.8182
$ws.Range("H136").Value = 1364.7894
$ws.Range("I136").Value = 1115.875
$ws.Range("J136").Value = 1545.8182
$ws.Range("K136").Value = 3347.625
$ws.Range("L136").Value = 4637.4546
$ws.Range("M136").Value = -797.625
$ws.Range("N136").Value = -9737.454600000001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1397.7858
$ws.Range("I5").Value = 659.4
$ws.Range("K5").Value = 1978.2
$ws.Range("M5").Value = -1866.2
$ws.Range("H92").Value = 537.75
$ws.Range("I92").Value = 367.33334
$ws.Range("J92").Value = 640
$ws.Range("K92").Value = 1102.00002
$ws.Range("L92").Value = 1920
$ws.Range("M92").Value = 145.9999800000001
$ws.Range("N92").Value = -4416
$ws.Range("H113").Value = 971.2826
$ws.Range("I113").Value = 699.8333
$ws.Range("J113").Value = 1480.25
$ws.Range("K113").Value = 2099.4999
$ws.Range("L113").Value = 4440.75
$ws.Range("M113").Value = 70.5001000000002
$ws.Range("N113").Value = -8780.75
$ws.Range("H120").Value = 17200
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 17200
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 51600
$ws.Range("M120").ClearContents()
$ws.Range("N120").Value = -61276
$ws.Range("H132").Value = 3600.5398
$ws.Range("I132").Value = 2711.3044
$ws.Range("K132").Value = 24401.7396
$ws.Range("M132").Value = -21871.7396
$ws.Range("H135").Value = 1397.7858
$ws.Range("I135").Value = 659.4
$ws.Range("K135").Value = 5934.599999999999
$ws.Range("M135").Value = -3399.599999999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4106.0513
$ws.Range("I122").Value = 2769.3333
$ws.Range("J122").Value = 5251.8096
$ws.Range("K122").Value = 8307.999899999999
$ws.Range("L122").Value = 15755.4288
$ws.Range("M122").Value = -5857.999899999999
$ws.Range("N122").Value = -20655.4288
$ws.Range("H126").Value = 2741.2
$ws.Range("I126").Value = 2816
$ws.Range("J126").Value = 2566.6667
$ws.Range("K126").Value = 8448
$ws.Range("L126").Value = 7700.000100000001
$ws.Range("M126").Value = -5978
$ws.Range("N126").Value = -12640.0001
$ws.Range("H132").Value = 2194.7
$ws.Range("I132").Value = 1606
$ws.Range("J132").Value = 3288
$ws.Range("K132").Value = 4818
$ws.Range("L132").Value = 9864
$ws.Range("M132").Value = -2288
$ws.Range("N132").Value = -14924
$ws.Range("H136").Value = 40326
$ws.Range("J136").Value = 40326
$ws.Range("L136").Value = 120978
$ws.Range("N136").Value = -126078

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 47622200
$ws.Range("I7").Value = 55558652
$ws.Range("J7").Value = 3470
$ws.Range("K7").Value = 55558652
$ws.Range("L7").Value = 3470
$ws.Range("M7").Value = -55558540
$ws.Range("N7").Value = -3694
$ws.Range("H40").Value = 55558320
$ws.Range("I40").Value = 71430550
$ws.Range("J40").Value = 5500
$ws.Range("K40").Value = 71430550
$ws.Range("L40").Value = 5500
$ws.Range("M40").Value = -71430414
$ws.Range("N40").Value = -5772
$ws.Range("H55").Value = 782.1739
$ws.Range("I55").Value = 170.3
$ws.Range("K55").Value = 170.3
$ws.Range("M55").Value = 2.699999999999989
$ws.Range("H82").Value = 2500
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 2500
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 2500
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = -3222
$ws.Range("H85").Value = 2500
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 2500
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 2500
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -4996
$ws.Range("H126").Value = 47622200
$ws.Range("I126").Value = 55558652
$ws.Range("J126").Value = 3470
$ws.Range("K126").Value = 166675956
$ws.Range("L126").Value = 10410
$ws.Range("M126").Value = -166673486
$ws.Range("N126").Value = -15350
$ws.Range("H136").Value = 5209838
$ws.Range("I136").Value = 1445.4
$ws.Range("J136").Value = 13890493
$ws.Range("K136").Value = 4336.200000000001
$ws.Range("L136").Value = 41671479
$ws.Range("M136").Value = -1786.200000000001
$ws.Range("N136").Value = -41676579

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 397.36365
$ws.Range("I113").Value = 413
$ws.Range("K113").Value = 1239
$ws.Range("M113").Value = 931
$ws.Range("H115").Value = 39600
$ws.Range("J115").Value = 39600
$ws.Range("L115").Value = 39600
$ws.Range("N115").Value = -42734
$ws.Range("H126").Value = 1391.5454
$ws.Range("I126").Value = 1101
$ws.Range("J126").Value = 1900
$ws.Range("K126").Value = 3303
$ws.Range("L126").Value = 5700
$ws.Range("M126").Value = -833
$ws.Range("N126").Value = -10640
$ws.Range("H136").Value = 5287.615
$ws.Range("I136").Value = 5305.4443
$ws.Range("K136").Value = 15916.3329
$ws.Range("M136").Value = -13366.3329
